# Bugfixed the naive forecaster component module:
# the forecast vector table is regenerated with corrected AR(2) coefficients
# and one fewer trailing observation (the data now ends at row 52 instead of 53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 53) - the data set shrank by one row
$ws.Rows("53:53").Delete()

# Clear cells that should no longer contain a value
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()

# Update cell values to match the new data
$ws.Range("A2").Value = 39583
$ws.Range("B2").Value = 2008
$ws.Range("D2").Value = 2009
$ws.Range("A3").Value = 39765
$ws.Range("A4").Value = 39948
$ws.Range("B4").Value = 2009
$ws.Range("D4").Value = 2010
$ws.Range("A5").Value = 40130
$ws.Range("C5").Value = -1.611885206309638
$ws.Range("A6").Value = 40310
$ws.Range("B6").Value = 2010
$ws.Range("D6").Value = 2011
$ws.Range("A7").Value = 40494
$ws.Range("A8").Value = 40676
$ws.Range("B8").Value = 2011
$ws.Range("D8").Value = 2012
$ws.Range("A9").Value = 40862
$ws.Range("A10").Value = 41044
$ws.Range("B10").Value = 2012
$ws.Range("D10").Value = 2013
$ws.Range("A11").Value = 41228
$ws.Range("E11").Value = 0.4163953164477929
$ws.Range("A12").Value = 41409
$ws.Range("B12").Value = 2013
$ws.Range("C12").Value = 0.03393100538855442
$ws.Range("D12").Value = 2014
$ws.Range("E12").Value = 0.3351240474928963
$ws.Range("A13").Value = 41592
$ws.Range("E13").Value = 1.653207170606596
$ws.Range("A14").Value = 41774
$ws.Range("B14").Value = 2014
$ws.Range("C14").Value = 1.743169463154315
$ws.Range("D14").Value = 2015
$ws.Range("E14").Value = 1.13752746419209
$ws.Range("A15").Value = 41957
$ws.Range("E15").Value = 1.270027657109818
$ws.Range("A16").Value = 42137
$ws.Range("B16").Value = 2015
$ws.Range("C16").Value = 1.311489985227077
$ws.Range("D16").Value = 2016
$ws.Range("E16").Value = 1.069982194174801
$ws.Range("A17").Value = 42321
$ws.Range("E17").Value = 1.579162878845075
$ws.Range("A18").Value = 42503
$ws.Range("B18").Value = 2016
$ws.Range("C18").Value = 1.721454720714122
$ws.Range("D18").Value = 2017
$ws.Range("E18").Value = 1.335637690776181
$ws.Range("A19").Value = 42689
$ws.Range("E19").Value = 1.452243308058287
$ws.Range("A20").Value = 42867
$ws.Range("B20").Value = 2017
$ws.Range("C20").Value = 1.818507532114921
$ws.Range("D20").Value = 2018
$ws.Range("E20").Value = 1.41495314213913
$ws.Range("A21").Value = 43053
$ws.Range("E21").Value = 2.068578555939404
$ws.Range("A22").Value = 43145
$ws.Range("B22").Value = 2018
$ws.Range("C22").Value = 2.23057583006443
$ws.Range("D22").Value = 2019
$ws.Range("E22").Value = 1.546132847114134
$ws.Range("A23").Value = 43235
$ws.Range("C23").Value = 2.466427116525516
$ws.Range("E23").Value = 1.859723853307749
$ws.Range("A24").Value = 43326
$ws.Range("C24").Value = 2.313955758667841
$ws.Range("E24").Value = 1.449675877460654
$ws.Range("A25").Value = 43418
$ws.Range("E25").Value = 1.651658474923545
$ws.Range("A26").Value = 43510
$ws.Range("B26").Value = 2019
$ws.Range("C26").Value = 1.34715816715496
$ws.Range("D26").Value = 2020
$ws.Range("E26").Value = 1.186882640643594
$ws.Range("A27").Value = 43600
$ws.Range("C27").Value = 1.128030950601477
$ws.Range("E27").Value = 0.9160236606447159
$ws.Range("A28").Value = 43691
$ws.Range("C28").Value = 0.9462474687678801
$ws.Range("E28").Value = 0.4351554058081408
$ws.Range("A29").Value = 43783
$ws.Range("E29").Value = 0.115841687688345
$ws.Range("A30").Value = 43875
$ws.Range("B30").Value = 2020
$ws.Range("C30").Value = 0.1938172373549873
$ws.Range("D30").Value = 2021
$ws.Range("E30").Value = 0.8212189468394859
$ws.Range("A31").Value = 43966
$ws.Range("C31").Value = 0.1798886261929367
$ws.Range("E31").Value = 0.8029144802146782
$ws.Range("A32").Value = 44068
$ws.Range("C32").Value = -4.43626840667447
$ws.Range("E32").Value = -14.53740902633983
$ws.Range("A33").Value = 44159
$ws.Range("E33").Value = -2.092304328310923
$ws.Range("A34").Value = 44251
$ws.Range("B34").Value = 2021
$ws.Range("C34").Value = -2.180664970010993
$ws.Range("D34").Value = 2022
$ws.Range("E34").Value = 0.6409010839486307
$ws.Range("A35").Value = 44341
$ws.Range("C35").Value = -1.929204335549095
$ws.Range("E35").Value = 0.8768515943972544
$ws.Range("A36").Value = 44432
$ws.Range("C36").Value = -1.513408827666285
$ws.Range("E36").Value = 1.518684466917741
$ws.Range("A37").Value = 44525
$ws.Range("E37").Value = 1.533339625605379
$ws.Range("A38").Value = 44617
$ws.Range("B38").Value = 2022
$ws.Range("C38").Value = 1.586146963184465
$ws.Range("D38").Value = 2023
$ws.Range("E38").Value = 0.9586871495637528
$ws.Range("A39").Value = 44706
$ws.Range("C39").Value = 1.632302710072264
$ws.Range("E39").Value = 0.9471575920676267
$ws.Range("A40").Value = 44798
$ws.Range("C40").Value = 1.618732201786743
$ws.Range("E40").Value = 0.8989718345186803
$ws.Range("A41").Value = 44890
$ws.Range("E41").Value = 0.492911192428136
$ws.Range("A42").Value = 44981
$ws.Range("B42").Value = 2023
$ws.Range("C42").Value = 0.1380617204474799
$ws.Range("D42").Value = 2024
$ws.Range("E42").Value = 0.6801419395370711
$ws.Range("A43").Value = 45071
$ws.Range("C43").Value = 0.01243672673012508
$ws.Range("E43").Value = 0.5779606211723021
$ws.Range("A44").Value = 45163
$ws.Range("C44").Value = -0.09609276733164585
$ws.Range("E44").Value = 0.3995394213445191
$ws.Range("A45").Value = 45254
$ws.Range("E45").Value = 0.2100922168233987
$ws.Range("A46").Value = 45345
$ws.Range("B46").Value = 2024
$ws.Range("C46").Value = 0.026532539029267
$ws.Range("D46").Value = 2025
$ws.Range("E46").Value = 0.5801787874785802
$ws.Range("A47").Value = 45436
$ws.Range("C47").Value = 0.01966607787367014
$ws.Range("E47").Value = 0.5811853063761419
$ws.Range("A48").Value = 45534
$ws.Range("C48").Value = -0.02761034355766023
$ws.Range("E48").Value = 0.5145439483255743
$ws.Range("A49").Value = 45618
$ws.Range("E49").Value = 0.5208382580577098
$ws.Range("A50").Value = 45713
$ws.Range("B50").Value = 2025
$ws.Range("C50").Value = 0.563860530038518
$ws.Range("D50").Value = 2026
$ws.Range("E50").Value = 0.6380016822143952
$ws.Range("A51").Value = 45800
$ws.Range("C51").Value = 0.5364374648222148
$ws.Range("E51").Value = 0.591050555601802
$ws.Range("A52").Value = 45891
$ws.Range("C52").Value = 0.5152269879013183
$ws.Range("E52").Value = 0.5355893905819142
